$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C25").Value = "IMAGE_RETRIEVAL_BY_METADATA"
$ws.Range("C26").Value = "IMAGE_RETRIEVAL_BY_METADATA"
$ws.Range("C47").Value = "BINARY_VISUAL_QA"
$ws.Range("C51").Value = "BINARY_VISUAL_QA"
$ws.Range("C54").Value = "BINARY_VISUAL_QA"
$ws.Range("C55").Value = "BINARY_VISUAL_QA"
$ws.Range("C56").Value = "BINARY_VISUAL_QA"
$ws.Range("C58").Value = "BINARY_VISUAL_QA"
$ws.Range("C59").Value = "BINARY_VISUAL_QA"
